$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data for A414:R445 (Apio / Vega Monumental Concepcion block)
# Each inner array is one row, columns A..R in order.
$rows = @(
  @(11, "Vega Monumental Concepción", "Bíobío", 45021, 8, 100112017, "Apio", "Americana (o)", "Primera", 100, 7500, 8000, 7750, "`$/docena de matas", "Región de Coquimbo", 1292, 6, "Hortaliza"),
  @(11, "Vega Monumental Concepción", "Bíobío", 45021, 8, 100112017, "Apio", "Americana (o)", "Segunda", 50, 6500, 6500, 6500, "`$/docena de matas", "Región de Coquimbo", 1083, 6, "Hortaliza"),
  @(11, "Vega Monumental Concepción", "Bíobío", 44264, 8, 100112017, "Apio", "Americana (o)", "Primera", 100, 7000, 8000, 7500, "`$/docena de matas", "Región de Coquimbo", 1250, 6, "Hortaliza"),
  @(11, "Vega Monumental Concepción", "Bíobío", 44264, 8, 100112017, "Apio", "Americana (o)", "Segunda", 50, 6000, 6000, 6000, "`$/docena de matas", "Región de Coquimbo", 1000, 6, "Hortaliza"),
  @(11, "Vega Monumental Concepción", "Bíobío", 44316, 8, 100112017, "Apio", "Americana (o)", "Primera", 100, 7000, 7500, 7250, "`$/docena de matas", "Región de Coquimbo", 1208, 6, "Hortaliza"),
  @(11, "Vega Monumental Concepción", "Bíobío", 44901, 8, 100112017, "Apio", "Americana (o)", "Primera", 220, 9000, 10000, 9455, "`$/docena de matas", "Región de Coquimbo", 1576, 6, "Hortaliza"),
  @(11, "Vega Monumental Concepción", "Bíobío", 44474, 8, 100112017, "Apio", "Americana (o)", "Primera", 100, 7000, 7500, 7250, "`$/docena de matas", "Región de Coquimbo", 1208, 6, "Hortaliza"),
  @(11, "Vega Monumental Concepción", "Bíobío", 44474, 8, 100112017, "Apio", "Americana (o)", "Segunda", 50, 6500, 6500, 6500, "`$/docena de matas", "Región de Coquimbo", 1083, 6, "Hortaliza"),
  @(11, "Vega Monumental Concepción", "Bíobío", 45002, 8, 100112017, "Apio", "Americana (o)", "Primera", 250, 8000, 8500, 8300, "`$/docena de matas", "Región de Coquimbo", 1383, 6, "Hortaliza"),
  @(11, "Vega Monumental Concepción", "Bíobío", 45002, 8, 100112017, "Apio", "Americana (o)", "Segunda", 220, 6000, 6500, 6227, "`$/docena de matas", "Región de Coquimbo", 1038, 6, "Hortaliza"),
  @(11, "Vega Monumental Concepción", "Bíobío", 44705, 8, 100112017, "Apio", "Americana (o)", "Primera", 100, 7000, 7500, 7250, "`$/docena de matas", "Región de Coquimbo", 1208, 6, "Hortaliza"),
  @(11, "Vega Monumental Concepción", "Bíobío", 44441, 8, 100112017, "Apio", "Americana (o)", "Primera", 100, 7000, 7500, 7250, "`$/docena de matas", "Región de Coquimbo", 1208, 6, "Hortaliza"),
  @(11, "Vega Monumental Concepción", "Bíobío", 44441, 8, 100112017, "Apio", "Americana (o)", "Segunda", 50, 6500, 6500, 6500, "`$/docena de matas", "Región de Coquimbo", 1083, 6, "Hortaliza"),
  @(11, "Vega Monumental Concepción", "Bíobío", 44952, 8, 100112017, "Apio", "Americana (o)", "Primera", 100, 8000, 8500, 8250, "`$/docena de matas", "Región de Coquimbo", 1375, 6, "Hortaliza"),
  @(11, "Vega Monumental Concepción", "Bíobío", 44952, 8, 100112017, "Apio", "Americana (o)", "Segunda", 50, 7000, 7000, 7000, "`$/docena de matas", "Región de Coquimbo", 1167, 6, "Hortaliza"),
  @(11, "Vega Monumental Concepción", "Bíobío", 44727, 8, 100112017, "Apio", "Americana (o)", "Primera", 200, 6000, 6500, 6250, "`$/docena de matas", "Región de Coquimbo", 1042, 6, "Hortaliza"),
  @(11, "Vega Monumental Concepción", "Bíobío", 44727, 8, 100112017, "Apio", "Americana (o)", "Segunda", 150, 5000, 5000, 5000, "`$/docena de matas", "Región de Coquimbo", 833, 6, "Hortaliza"),
  @(11, "Vega Monumental Concepción", "Bíobío", 44252, 8, 100112017, "Apio", "Americana (o)", "Primera", 100, 7000, 8000, 7500, "`$/docena de matas", "Región de Coquimbo", 1250, 6, "Hortaliza"),
  @(11, "Vega Monumental Concepción", "Bíobío", 44252, 8, 100112017, "Apio", "Americana (o)", "Segunda", 50, 6000, 6000, 6000, "`$/docena de matas", "Región de Coquimbo", 1000, 6, "Hortaliza"),
  @(11, "Vega Monumental Concepción", "Bíobío", 44855, 8, 100112017, "Apio", "Americana (o)", "Primera", 100, 8000, 8500, 8250, "`$/docena de matas", "Región de Coquimbo", 1375, 6, "Hortaliza"),
  @(11, "Vega Monumental Concepción", "Bíobío", 44855, 8, 100112017, "Apio", "Americana (o)", "Segunda", 50, 7000, 7000, 7000, "`$/docena de matas", "Región de Coquimbo", 1167, 6, "Hortaliza"),
  @(11, "Vega Monumental Concepción", "Bíobío", 44391, 8, 100112017, "Apio", "Americana (o)", "Primera", 100, 8000, 9000, 8500, "`$/docena de matas", "Región de Coquimbo", 1417, 6, "Hortaliza"),
  @(11, "Vega Monumental Concepción", "Bíobío", 44391, 8, 100112017, "Apio", "Americana (o)", "Segunda", 50, 7000, 7000, 7000, "`$/docena de matas", "Provincia de Limarí", 1167, 6, "Hortaliza"),
  @(11, "Vega Monumental Concepción", "Bíobío", 44168, 8, 100112017, "Apio", "Americana (o)", "Primera", 100, 7500, 8000, 7750, "`$/docena de matas", "Provincia de Limarí", 1292, 6, "Hortaliza"),
  @(11, "Vega Monumental Concepción", "Bíobío", 44875, 8, 100112017, "Apio", "Americana (o)", "Primera", 270, 8000, 8500, 8222, "`$/docena de matas", "Región de Coquimbo", 1370, 6, "Hortaliza"),
  @(11, "Vega Monumental Concepción", "Bíobío", 44875, 8, 100112017, "Apio", "Americana (o)", "Segunda", 150, 7000, 7000, 7000, "`$/docena de matas", "Región de Coquimbo", 1167, 6, "Hortaliza"),
  @(11, "Vega Monumental Concepción", "Bíobío", 44642, 8, 100112017, "Apio", "Americana (o)", "Primera", 220, 7500, 8000, 7773, "`$/docena de matas", "Región de Coquimbo", 1296, 6, "Hortaliza"),
  @(11, "Vega Monumental Concepción", "Bíobío", 44642, 8, 100112017, "Apio", "Americana (o)", "Segunda", 200, 5500, 6500, 6000, "`$/docena de matas", "Región de Coquimbo", 1000, 6, "Hortaliza"),
  @(11, "Vega Monumental Concepción", "Bíobío", 44243, 8, 100112017, "Apio", "Americana (o)", "Primera", 100, 8000, 9000, 8500, "`$/docena de matas", "Región de Coquimbo", 1417, 6, "Hortaliza"),
  @(11, "Vega Monumental Concepción", "Bíobío", 44243, 8, 100112017, "Apio", "Americana (o)", "Segunda", 50, 7000, 7000, 7000, "`$/docena de matas", "Región de Coquimbo", 1167, 6, "Hortaliza"),
  @(11, "Vega Monumental Concepción", "Bíobío", 44217, 8, 100112017, "Apio", "Americana (o)", "Primera", 100, 8000, 9000, 8500, "`$/docena de matas", "Región de Coquimbo", 1417, 6, "Hortaliza"),
  @(11, "Vega Monumental Concepción", "Bíobío", 44217, 8, 100112017, "Apio", "Americana (o)", "Segunda", 50, 7000, 7000, 7000, "`$/docena de matas", "Región de Coquimbo", 1167, 6, "Hortaliza")
)

$startRow = 414
$r = $startRow
foreach ($row in $rows) {
  $c = 1
  foreach ($val in $row) {
    $ws.Cells.Item($r, $c).Value = $val
    $c = $c + 1
  }
  $r = $r + 1
}

# Preserve the date-time number format on column D for the two newly appended rows
# (rows 414:443 already carry this style on their D cell; 444:445 are brand new).
$ws.Range("D444:D445").NumberFormat = "YYYY-MM-DD HH:MM:SS"

